$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (shifts existing data right by 2)
$ws.Range("A1:B1").EntireColumn.Insert()

# Add new headers
$ws.Range("A1").Value = "Sno"
$ws.Range("B1").Value = "ExecutionFlag"

# Fill Sno column (1-5) for rows 2-6
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Fill ExecutionFlag column
$ws.Range("B2").Value = "Y"
$ws.Range("B3").Value = "Y"
$ws.Range("B4").Value = "N"
$ws.Range("B5").Value = "Y"
$ws.Range("B6").Value = "Y"

# Update selection to match target state
$ws.Range("B6").Select()
